$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append daily COVID data rows through 2021-09-09 (serials 44441-44448):
# col A = date serial, col B = "nuovi pos.", col C = "somma mobile 7gg.",
# col D = "somma mobile 7gg. per 100mila abitanti"
$data = @(
    @(367, 44441, 0, 5, 728.862973760933),
    @(368, 44442, 0, 4, 583.0903790087464),
    @(369, 44443, 0, 0, 0),
    @(370, 44444, 0, 0, 0),
    @(371, 44445, 0, 0, 0),
    @(372, 44446, 0, 0, 0),
    @(373, 44447, 0, 0, 0),
    @(374, 44448, 0, 0, 0)
)

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value = $r[1]
    $ws.Cells.Item($row, 2).Value = $r[2]
    $ws.Cells.Item($row, 3).Value = $r[3]
    $ws.Cells.Item($row, 4).Value = $r[4]
}

# Match the date-format/border/font styling used by the rest of column A
# (same style as the preceding row, index "s=2" in the sheet XML).
$ws.Range("A366").Copy()
$ws.Range("A367:A374").PasteSpecial(-4122)

$wb.Save()
